$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Issue list")
$lo = $ws.ListObjects.Item("Table1")

# Add a new row to the table (expands Table1/autoFilter ref from A1:M39 to
# A1:M40, same as Excel does when a ListRow is appended).
$lo.ListRows.Add() | Out-Null

# Copy the formatting of the previous last row down into the new row so the
# per-cell styles (alignment etc.) match the rest of the table, then
# overwrite with the real values below.
$ws.Range("A39:M39").Copy()
$ws.Range("A40:M40").PasteSpecial(-4122)

$ws.Range("A40").Value = 39
$ws.Range("B40").Value = 5493
$ws.Range("C40").Value = "Bug 5493 - Extend data_environment for dimensional_characteristic_representation and others"
$ws.Range("C40").HorizontalAlignment = -4131
$ws.Range("C40").VerticalAlignment = -4160
$ws.Range("C40").WrapText = $false
$ws.Range("D40").Value = "valid-shtolo_issue"
$ws.Range("E40").Value = "x"
$ws.Range("F40").Value = 1.49
$ws.Range("G40").Value = "x"
$ws.Range("H40").Value = 1.71
$ws.Range("I40").Value = "x"
$ws.Range("J40").Value = 1.58
$ws.Range("K40").Value = "change data_environment.elements as :`nENTITY data_environment;`n  name : label;`n  description : text;`n  elements : SET [1:?] OF property_definition_representation;`nEND_ENTITY;"
